# Apply the "add fao to total luc differences" edit.
#
# 1) Rename the three bookkeeping-model labels (and the "Bookkeeping average"
#    label) in the "data" sheet, column B, to be prefixed with "Bookkeeping | ".
# 2) Append 33 new rows (years 1990-2022) for a new "FAO" series, using the
#    same units ("GtCO2/year") as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Step 1: rename the existing variable labels ------------------------
$lastRow = 166

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 2)
    $b = $cell.Text
    if ($b -eq "BLUE") {
        $cell.Value = "Bookkeeping | BLUE"
    } elseif ($b -eq "H&N") {
        $cell.Value = "Bookkeeping | H&N"
    } elseif ($b -eq "OSCAR") {
        $cell.Value = "Bookkeeping | OSCAR"
    } elseif ($b -eq "Bookkeeping average") {
        $cell.Value = "Bookkeeping | average"
    }
}

# --- Step 2: append the new FAO rows -------------------------------------
$years = @(1990,1991,1992,1993,1994,1995,1996,1997,1998,1999,2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022)
$values = @(1.7635435932,1.7635435932,1.7635435931,1.7639772221,1.7636441844,1.7692517676,1.6361572212,2.3259997112,1.7700120542,1.6736788596,1.6215367559,1.3444369577,1.7201421788,1.4476078186,1.7178166743,1.5112234373,1.8068813722,1.3808307486,1.3687807169,1.6726327904,1.3818374722,0.3066388384,0.3161867504,0.3019906273,0.5654136678,0.5815682263,1.1828981771,1.1479663139,1.2958760425,1.4706543365,1.1731730705,1.1535669849,1.145102449)

$startRow = $lastRow + 1

for ($i = 0; $i -lt $years.Count; $i++) {
    $r = $startRow + $i
    $data.Cells.Item($r, 1).Value = $years[$i]
    $data.Cells.Item($r, 2).Value = "FAO"
    $data.Cells.Item($r, 3).Value = $values[$i]
    $data.Cells.Item($r, 4).Value = "GtCO2/year"
}
